$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = -6
$ws.Range("F4").Value = -1
$ws.Range("F5").Value = 2
$ws.Range("F6").Value = 2
$ws.Range("F11").Value = -5
$ws.Range("F14").Value = 2
